$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "x"
$ws.Range("B14").Value = "x"
$ws.Range("C14").Value = "Nick"

$ws.Range("C14").Select()
